# Apply weekly fruit/vegetable price update (Caqui) to rows 2, 4 and 5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: now carries the figures previously held by row 5
$ws.Range("D2").Value = 44355
$ws.Range("M2").Value = 270
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 1139

# Row 4: date updated, quality changes from Primera to Segunda
$ws.Range("D4").Value = 44305
$ws.Range("L4").Value = "Segunda"

# Row 5: now carries the figures previously held by row 4
$ws.Range("D5").Value = 44342
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 24500
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1361
